$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.577.20'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.60%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.599.79'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.26%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '609.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.60%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.18%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.490'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E9').Value = '  +0.34%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '8.05'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.41%  '

$ws.Range('E11').Value = '  +1.04%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.203.51'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.00%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000210'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.94%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '30.08'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.20%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.600.26'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.13%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.650.66'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.52%  '

$ws.Range('E17').Value = '  +0.81%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.52'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.50%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.34'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.09%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.08'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.67%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '429.47'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.24%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.622'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.26%  '

$ws.Range('E23').Value = '  +0.61%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.735.01'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.92%  '

$ws.Range('E26').Value = '  +2.59%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.31'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.44%  '

$ws.Range('E28').Value = '  +1.21%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.05%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.52%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.593.26'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.11%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.57'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.70%  '

$ws.Range('E33').Value = '  -1.72%  '

$ws.Range('E34').Value = '  -3.46%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.86'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.41%  '

$ws.Range('E36').Value = '  +0.03%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.73'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.66%  '

$ws.Range('E38').Value = '  -0.12%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '177.11'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.50%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0860'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.44%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.25'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.99%  '

$ws.Range('E42').Value = '  +0.24%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.92'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.65%  '

$ws.Range('E44').Value = '  +8.57%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.06%  '

$ws.Range('E46').Value = '  -1.78%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '25.10'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.01%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '24.26'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.61%  '

$ws.Range('E49').Value = '  +1.24%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.955'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.02%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.237'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.03%  '
